# Commit: "remove outliers in exposure to coast"
# The three provinces Iloilo, Sulu, and Tawi-Tawi are dropped from the
# results table (they were exposure-to-coast outliers). Removing their
# rows shifts every following row up by one, which is why many of the
# later "Risk to well-being" figures shift slightly after recalculation
# upstream; the structural change we need to reproduce here is simply
# deleting those three whole rows from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$provincesToRemove = @("Iloilo", "Sulu", "Tawi-Tawi")

# Column A holds the province name; column A1 is the header ("province").
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Collect the row numbers that must be removed.
$rowsToDelete = New-Object System.Collections.ArrayList
for ($r = 1; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($provincesToRemove -contains $name) {
        [void]$rowsToDelete.Add($r)
    }
}

# Delete from the bottom up so earlier row numbers stay valid.
$sortedRows = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sortedRows) {
    $ws.Rows.Item($r).Delete()
}
